$d = $word.ActiveDocument

$d.Content.Find.Execute("85×44=3740", $true, $true, $false, $false, $false, $true, 1, $false, "81×60=4860", 2) | Out-Null
$d.Content.Find.Execute("39×99=3861", $true, $true, $false, $false, $false, $true, 1, $false, "29×36=1044", 2) | Out-Null
$d.Content.Find.Execute("90×50=4500", $true, $true, $false, $false, $false, $true, 1, $false, "91×73=6643", 2) | Out-Null
$d.Content.Find.Execute("58×23=1334", $true, $true, $false, $false, $false, $true, 1, $false, "37×47=1739", 2) | Out-Null
$d.Content.Find.Execute("62×31=1922", $true, $true, $false, $false, $false, $true, 1, $false, "53×46=2438", 2) | Out-Null
$d.Content.Find.Execute("26×23=598", $true, $true, $false, $false, $false, $true, 1, $false, "79×90=7110", 2) | Out-Null
$d.Content.Find.Execute("44×98=4312", $true, $true, $false, $false, $false, $true, 1, $false, "29×94=2726", 2) | Out-Null
$d.Content.Find.Execute("71×46=3266", $true, $true, $false, $false, $false, $true, 1, $false, "92×86=7912", 2) | Out-Null
$d.Content.Find.Execute("27×74=1998", $true, $true, $false, $false, $false, $true, 1, $false, "84×83=6972", 2) | Out-Null
$d.Content.Find.Execute("73×32=2336", $true, $true, $false, $false, $false, $true, 1, $false, "76×61=4636", 2) | Out-Null
$d.Content.Find.Execute("59×28=1652", $true, $true, $false, $false, $false, $true, 1, $false, "81×77=6237", 2) | Out-Null
$d.Content.Find.Execute("89×72=6408", $true, $true, $false, $false, $false, $true, 1, $false, "44×90=3960", 2) | Out-Null
$d.Content.Find.Execute("92×98=9016", $true, $true, $false, $false, $false, $true, 1, $false, "40×35=1400", 2) | Out-Null
$d.Content.Find.Execute("73×96=7008", $true, $true, $false, $false, $false, $true, 1, $false, "34×95=3230", 2) | Out-Null
$d.Content.Find.Execute("39×30=1170", $true, $true, $false, $false, $false, $true, 1, $false, "21×21=441", 2) | Out-Null
$d.Content.Find.Execute("64×39=2496", $true, $true, $false, $false, $false, $true, 1, $false, "50×16=800", 2) | Out-Null
$d.Content.Find.Execute("65×79=5135", $true, $true, $false, $false, $false, $true, 1, $false, "20×23=460", 2) | Out-Null
$d.Content.Find.Execute("88×20=1760", $true, $true, $false, $false, $false, $true, 1, $false, "70×81=5670", 2) | Out-Null
$d.Content.Find.Execute("52×57=2964", $true, $true, $false, $false, $false, $true, 1, $false, "31×12=372", 2) | Out-Null
$d.Content.Find.Execute("81×23=1863", $true, $true, $false, $false, $false, $true, 1, $false, "83×22=1826", 2) | Out-Null
$d.Content.Find.Execute("98×28=2744", $true, $true, $false, $false, $false, $true, 1, $false, "69×96=6624", 2) | Out-Null
$d.Content.Find.Execute("29×96=2784", $true, $true, $false, $false, $false, $true, 1, $false, "31×88=2728", 2) | Out-Null
$d.Content.Find.Execute("30×32=960", $true, $true, $false, $false, $false, $true, 1, $false, "30×85=2550", 2) | Out-Null
$d.Content.Find.Execute("39×92=3588", $true, $true, $false, $false, $false, $true, 1, $false, "77×38=2926", 2) | Out-Null
$d.Content.Find.Execute("53×75=3975", $true, $true, $false, $false, $false, $true, 1, $false, "72×92=6624", 2) | Out-Null
